$wb = $excel.ActiveWorkbook

# ---- Week 2: grade Tutorial_02_1/2/3 and Program_02_1/2/3 (all pass) ----
$ws2 = $wb.Worksheets.Item("Week 2")

$ws2.Range("B2").Value = 1
$ws2.Range("D2").Style = "Normal"

$ws2.Range("B3").Value = 1
$ws2.Range("D3").Style = "Normal"

$ws2.Range("B4").Value = 1
$ws2.Range("D4").Style = "Normal"

$ws2.Range("B5").Value = 1
$ws2.Range("D5").Style = "Normal"

$ws2.Range("B6").Value = 1
$ws2.Range("D6").Style = "Normal"

$ws2.Range("B7").Value = 1
$ws2.Range("D7").Style = "Normal"

$ws2.Columns.Item(2).ColumnWidth = 8.3
$ws2.Columns.Item(4).ColumnWidth = 5.8

# ---- Week 3: grade Tutorial_03_1-5 and Program_03_1-6 ----
$ws3 = $wb.Worksheets.Item("Week 3")

# Tutorial_03_1 .. Tutorial_03_5 (rows 2-6) -- pass
$ws3.Range("B2").Value = 1
$ws3.Range("D2").Style = "Normal"

$ws3.Range("B3").Value = 1
$ws3.Range("D3").Style = "Normal"

$ws3.Range("B4").Value = 1
$ws3.Range("D4").Style = "Normal"

$ws3.Range("B5").Value = 1
$ws3.Range("D5").Style = "Normal"

$ws3.Range("B6").Value = 1
$ws3.Range("D6").Style = "Normal"

# Program_03_1 (row 7) -- pass
$ws3.Range("B7").Value = 1
$ws3.Range("D7").Style = "Normal"

# Program_03_2 .. Program_03_6 (rows 8-12) -- fail, not found
$ws3.Range("B8").Value = 0
$ws3.Range("D8").Value = "Not Found"

$ws3.Range("B9").Value = 0
$ws3.Range("D9").Value = "Not Found"

$ws3.Range("B10").Value = 0
$ws3.Range("D10").Value = "Not Found"

$ws3.Range("B11").Value = 0
$ws3.Range("D11").Value = "Not Found"

$ws3.Range("B12").Value = 0
$ws3.Range("D12").Value = "Not Found"

$ws3.Columns.Item(2).ColumnWidth = 8.3
$ws3.Columns.Item(4).ColumnWidth = 9.8
